$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.813.90"
$ws.Range("E2").Value = "  +1.93%  "
$ws.Range("D3").Value = "2.117.08"
$ws.Range("E3").Value = "  +6.42%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5335"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.67%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4412"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09020"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "47.52"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +11.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.181"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.01"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.44%  "
$ws.Range("D13").Value = "2.110.44"
$ws.Range("E13").Value = "  +5.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.775"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.833"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "96.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001134"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06683"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.56%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("E22").Value = "  +4.52%  "
$ws.Range("D23").Value = "30.874.86"
$ws.Range("E23").Value = "  +1.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.69%  "
$ws.Range("D25").Value = "2.360.05"
$ws.Range("E25").Value = "  +6.06%  "
$ws.Range("E26").Value = "  +3.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.594"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +10.13%  "
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.55"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.188"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1090"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.232"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.013"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.561"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +19.12%  "
$ws.Range("E36").Value = "  +5.06%  "
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "12.92"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.87%  "
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.550"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.16%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06776"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.14%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.571"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2316"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6859"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.246"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6466"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.90%  "
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.0000"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "14.07"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.237"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.671"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.269"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.76%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "83.04"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.197"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +9.60%  "
